$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1748251748251748
$ws.Range("C2").Value = 0.6118881118881119
$ws.Range("J2").Value = 0.02097902097902098
$ws.Range("P2").Value = 0.1223776223776224
$ws.Range("S2").Value = 0.06993006993006994
$ws.Range("B3").Value = 0.00558659217877095
$ws.Range("C3").Value = 0.01675977653631285
$ws.Range("J3").Value = 0.0223463687150838
$ws.Range("P3").Value = 0.7430167597765364
$ws.Range("S3").Value = 0.2122905027932961
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("P4").Value = 0.6764705882352942
$ws.Range("S4").Value = 0.2647058823529412
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.07239819004524888
$ws.Range("D6").Value = 0.01809954751131222
$ws.Range("F6").Value = 0.03619909502262444
$ws.Range("J6").Value = 0.2443438914027149
$ws.Range("O6").Value = 0.01809954751131222
$ws.Range("Q6").Value = 0.167420814479638
$ws.Range("R6").Value = 0.04977375565610859
$ws.Range("S6").Value = 0.3936651583710407
$ws.Range("B7").Value = 0.08372093023255814
$ws.Range("D7").Value = 0.004651162790697674
$ws.Range("F7").Value = 0.05581395348837209
$ws.Range("J7").Value = 0.1441860465116279
$ws.Range("O7").Value = 0.03255813953488372
$ws.Range("Q7").Value = 0.1767441860465116
$ws.Range("R7").Value = 0.06976744186046512
$ws.Range("S7").Value = 0.4325581395348837
$ws.Range("B8").Value = 0.08089887640449438
$ws.Range("D8").Value = 0.01123595505617977
$ws.Range("F8").Value = 0.04943820224719101
$ws.Range("J8").Value = 0.1460674157303371
$ws.Range("O8").Value = 0.02247191011235955
$ws.Range("Q8").Value = 0.2224719101123596
$ws.Range("R8").Value = 0.07865168539325842
$ws.Range("S8").Value = 0.3887640449438202
$ws.Range("B9").Value = 0.0594059405940594
$ws.Range("D9").Value = 0.01485148514851485
$ws.Range("F9").Value = 0.05445544554455446
$ws.Range("J9").Value = 0.1089108910891089
$ws.Range("O9").Value = 0.02475247524752475
$ws.Range("Q9").Value = 0.1584158415841584
$ws.Range("R9").Value = 0.08415841584158416
$ws.Range("S9").Value = 0.495049504950495
$ws.Range("B10").Value = 0.09853528628495339
$ws.Range("D10").Value = 0.01464713715046605
$ws.Range("E10").Value = 0.001997336884154461
$ws.Range("F10").Value = 0.06724367509986684
$ws.Range("J10").Value = 0.1225033288948069
$ws.Range("O10").Value = 0.01198402130492676
$ws.Range("Q10").Value = 0.2463382157123835
$ws.Range("R10").Value = 0.05725699067909454
$ws.Range("S10").Value = 0.3794940079893475
$ws.Range("F11").Value = 0.00303030303030303
$ws.Range("G11").Value = 0.1303030303030303
$ws.Range("J11").Value = 0.08181818181818182
$ws.Range("K11").Value = 0.1848484848484848
$ws.Range("L11").Value = 0.5787878787878787
$ws.Range("S11").Value = 0.02121212121212121
$ws.Range("G12").Value = 0.7305699481865285
$ws.Range("J12").Value = 0.227979274611399
$ws.Range("L12").Value = 0.0155440414507772
$ws.Range("S12").Value = 0.02590673575129534
$ws.Range("F15").Value = 0.0158102766798419
$ws.Range("H15").Value = 0.1383399209486166
$ws.Range("I15").Value = 0.05928853754940711
$ws.Range("J15").Value = 0.3517786561264822
$ws.Range("K15").Value = 0.07509881422924901
$ws.Range("O15").Value = 0.08300395256916997
$ws.Range("S15").Value = 0.2766798418972332
$ws.Range("F16").Value = 0.01639344262295082
$ws.Range("H16").Value = 0.1693989071038251
$ws.Range("I16").Value = 0.1092896174863388
$ws.Range("J16").Value = 0.4480874316939891
$ws.Range("K16").Value = 0.09289617486338798
$ws.Range("M16").Value = 0.03825136612021858
$ws.Range("O16").Value = 0.03825136612021858
$ws.Range("S16").Value = 0.08743169398907104
$ws.Range("F17").Value = 0.02434782608695652
$ws.Range("H17").Value = 0.1686956521739131
$ws.Range("I17").Value = 0.08173913043478261
$ws.Range("J17").Value = 0.4365217391304348
$ws.Range("K17").Value = 0.1147826086956522
$ws.Range("M17").Value = 0.01391304347826087
$ws.Range("O17").Value = 0.07130434782608695
$ws.Range("S17").Value = 0.08869565217391304
$ws.Range("F18").Value = 0.01219512195121951
$ws.Range("H18").Value = 0.1829268292682927
$ws.Range("I18").Value = 0.1097560975609756
$ws.Range("J18").Value = 0.4634146341463415
$ws.Range("K18").Value = 0.1219512195121951
$ws.Range("M18").Value = 0.01219512195121951
$ws.Range("O18").Value = 0.0426829268292683
$ws.Range("S18").Value = 0.05487804878048781
$ws.Range("F19").Value = 0.01228323699421965
$ws.Range("H19").Value = 0.1871387283236994
$ws.Range("I19").Value = 0.07586705202312138
$ws.Range("J19").Value = 0.4082369942196532
$ws.Range("K19").Value = 0.1047687861271676
$ws.Range("M19").Value = 0.02601156069364162
$ws.Range("N19").Value = 0.002167630057803468
$ws.Range("O19").Value = 0.07369942196531792
$ws.Range("S19").Value = 0.1098265895953757
